$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Comentarios" header (F3) to "Notas"
$ws.Range("F3").Value = "Notas"

# Update the active selection to match the edited workbook's saved view
$ws.Range("F4").Select()
